# "excel import basic feature added"
# Adds a second sheet of imported evaluators and appends a batch of newly
# imported freelancer/employee rows to the existing Sheet1 list.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Sheet1: append newly imported rows 5-16
# ---------------------------------------------------------------------

# Row 5 - Nataliya (full record, own company)
$ws1.Range("B5").Value = "Nataliya"
$ws1.Range("E5").Value = "Nata INC"
$ws1.Range("A5").Value = "nata@fontanille.com"
$ws1.Range("C5").Value = "Fontanille"
$ws1.Range("D5").Value = "CEO"
$ws1.Range("F5").Value = "Management"

# Rows 6-16 - freelancer roster, first/last names entered per row
$ws1.Range("B6").Value = "Nicolas"
$ws1.Range("C6").Value = "Aussenac"
$ws1.Range("B7").Value = "Jean-Christophe"
$ws1.Range("C7").Value = "Bouvier"
$ws1.Range("B8").Value = "Manuel "
$ws1.Range("C8").Value = "Martinez"
$ws1.Range("B9").Value = "Aurelie "
$ws1.Range("C9").Value = "Lejeune"
$ws1.Range("B10").Value = "Yoann"
$ws1.Range("C10").Value = "Fontanille"
$ws1.Range("C11").Value = "Garnier"
$ws1.Range("B12").Value = "Carl"
$ws1.Range("C12").Value = "Cox"
$ws1.Range("B13").Value = "David"
$ws1.Range("C13").Value = "Guetta"
$ws1.Range("B14").Value = "Nora"
$ws1.Range("C14").Value = "En Pure"
$ws1.Range("B15").Value = "The"
$ws1.Range("C15").Value = "Prodigy"
$ws1.Range("B16").Value = "Chemical"
$ws1.Range("C16").Value = "Brothers"

# job_title/company filled as "Freelancer" for the whole batch (row 11 has
# no first name, row 12 has no job_title, row 15 has no hyperlinked email)
$ws1.Range("D6:D11").Value = "Freelancer"
$ws1.Range("D13:D16").Value = "Freelancer"
$ws1.Range("E6:E16").Value = "Freelancer"

# department filled as "None" for the whole batch
$ws1.Range("F6:F16").Value = "None"

# shared freelancer contact email, filled down (A15 intentionally left blank)
$ws1.Range("A6:A14").Value = "freelancer@gmail.com"
$ws1.Range("A16").Value = "freelancer@gmail.com"

# mailto hyperlinks for the two newly introduced addresses
$ws1.Hyperlinks.Add($ws1.Range("A5"), "mailto:nata@fontanille.com")
$ws1.Hyperlinks.Add($ws1.Range("A6"), "mailto:freelancer@gmail.com")
# Hyperlinks.Add() re-applies the Hyperlink style with a brand new cellXf;
# put it back on the shared "Hyperlink" style so A5/A6 line up with A2:A4.
$ws1.Range("A5:A6").Style = "Hyperlink"

# move the selection cursor the way it ends up after this kind of paste
$ws1.Range("A21").Select()

# ---------------------------------------------------------------------
# Sheet2: new sheet holding the rest of the imported records
# ---------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "email"
$ws2.Range("B1").Value = "first_name"
$ws2.Range("C1").Value = "last_name"
$ws2.Range("D1").Value = "job_title"
$ws2.Range("E1").Value = "company"
$ws2.Range("F1").Value = "department"
$ws2.Range("A1:F1").Style = "Normal"
$ws2.Range("A1:F1").Font.Bold = $true

$ws2.Range("A2").Value = "john@yies.co"
$ws2.Range("B2").Value = "John"
$ws2.Range("C2").Value = "Doo"
$ws2.Range("D2").Value = "CEO"
$ws2.Range("E2").Value = "YIES"
$ws2.Range("F2").Value = "Management"

$ws2.Range("A3").Value = "laura@ibm.fr"
$ws2.Range("B3").Value = "Laura"
$ws2.Range("C3").Value = "Williams"
$ws2.Range("D3").Value = "Engineer"
$ws2.Range("E3").Value = "IBM"
$ws2.Range("F3").Value = "Engineering"

$ws2.Range("A4").Value = "nicolas@nike.no"
$ws2.Range("B4").Value = "Nicolas"
$ws2.Range("C4").Value = "Dupont"
$ws2.Range("D4").Value = "Sales Manager"
$ws2.Range("E4").Value = "Nike"
$ws2.Range("F4").Value = "Sales"

$ws2.Columns.Item(1).ColumnWidth = 14.6640625
$ws2.Columns.Item(2).ColumnWidth = 10.88671875
$ws2.Columns.Item(3).ColumnWidth = 11
$ws2.Columns.Item(4).ColumnWidth = 12.77734375
$ws2.Columns.Item(5).ColumnWidth = 10
$ws2.Columns.Item(6).ColumnWidth = 11.77734375

$ws2.Range("C12").Select()
